$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: replace the sample/test data with a single draft input row ---
$ws.Range("A2").Value = "13088334935"
$ws.Range("B2").Value = "267794"
$ws.Range("C2").Value = "13088334935"
$ws.Range("D2").Value = "8 333,33"
$ws.Range("E2").Value = "'False"
$ws.Range("F2").Value = "7100"
$ws.Range("G2").Value = "1"
$ws.Range("H2").Value = "'False"
$ws.Range("I2").Value = "0"
$ws.Range("J2").Value = "0"
$ws.Range("K2").Value = "0"
$ws.Range("L2").Value = "FJELLHAMAR"
$ws.Range("P2").Value = "984661185"
$ws.Range("Q2").Value = "POSTEN NORGE AS"
$ws.Range("R2").Value = "Lørenskogveien 50`r"
$ws.Range("S2").Value = "1470 LØRENSKOG`r"
$ws.Range("T2").Clear()

# --- Rows 3-6: clear the remaining sample rows down to blank text cells ---
# (matching the already-blank rows 7:300) while preserving each column's
# original cell style (e.g. the quotePrefix style on columns E/H).
$ws.Range("A3:T6").Value = "'"
$ws.Range("A7:T7").Copy()
$ws.Range("A3:T6").PasteSpecial(-4122)

# --- Update the active selection shown when the sheet is opened ---
$ws.Range("A2:AB3").Select()
